$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 4).Value = -0.009519999999999999   # D
    $ws.Cells.Item($r, 5).Value = -0.0411                  # E

    $ws.Cells.Item($r, 7).Value = 0.0299210316597841       # G
    $ws.Cells.Item($r, 8).Value = 0.0299210316597841       # H
    $ws.Cells.Item($r, 9).Value = 0.03607911323625299      # I
    $ws.Cells.Item($r, 10).Value = 0.02749801928148581     # J
    $ws.Cells.Item($r, 11).Value = 128.4                   # K
    $ws.Cells.Item($r, 12).Value = 0.02325581395348838     # L
    $ws.Cells.Item($r, 13).Value = 134.0063                # M
    $ws.Cells.Item($r, 14).Value = 0.1611621166566446      # N
    $ws.Cells.Item($r, 15).Value = 1.04366277258567        # O
    $ws.Cells.Item($r, 16).Value = 49.2063                 # P
    $ws.Cells.Item($r, 17).Value = 0.05917775105231509     # Q
    $ws.Cells.Item($r, 18).Value = 0.3832266355140186      # R
    $ws.Cells.Item($r, 19).Value = 84.80000000000001       # S
    $ws.Cells.Item($r, 20).Value = 0.6328060695653861      # T
    $ws.Cells.Item($r, 21).Value = 0.002                   # U
    $ws.Cells.Item($r, 22).Value = 0.000002405291641611546 # V
    $ws.Cells.Item($r, 23).Value = 0.06408144931876029     # W
    $ws.Cells.Item($r, 24).Value = 0.07441555875720936     # X
    $ws.Cells.Item($r, 25).Value = -0.01033410943844906    # Y
    $ws.Cells.Item($r, 26).Value = 2.755505071123492       # Z
    $ws.Cells.Item($r, 27).Value = 0.07577093157598573     # AA
    $ws.Cells.Item($r, 28).Value = 0.07441555875720936     # AB
    $ws.Cells.Item($r, 29).Value = 0.001355372818776374    # AC

    $ws.Cells.Item($r, 36).Value = -0.000002405297427053342 # AJ
    $ws.Cells.Item($r, 37).Value = -0.0000009477334480722626 # AK
    $ws.Cells.Item($r, 38).Value = 0.344                   # AL
    $ws.Cells.Item($r, 39).Value = 0.344                   # AM

    $ws.Cells.Item($r, 41).Value = 579.0697674418604       # AO
    $ws.Cells.Item($r, 42).Value = -0.00000970873786407767 # AP
    $ws.Cells.Item($r, 43).Value = 579.0697674418604       # AQ
}
